$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "diet_test_rf"

$ws.Range("A1").Value = "Diet_orgs"
$ws.Range("B1").Value = "RF slurm"
$ws.Range("C1").Value = "d1_diet_orgs"

$ws.Range("A2").Value = "    d1d2_nutr_only_2015.csv,"
$ws.Range("B2").Value = "    rf_cardio_nut_only.slurm,"
$ws.Range("C2").Value = "    d1_nutr_only_2015.csv,"

$ws.Range("A3").Value = "    d1d2_food_g_2015.csv,"
$ws.Range("B3").Value = "    rf_cardio_food_g.slurm,"
$ws.Range("C3").Value = "    d1_food_g_2015.csv,"

$ws.Range("A4").Value = "    d1d2_nutri_food_g_2015.csv,"
$ws.Range("B4").Value = "    rf_cardio_nutr_food_g.slurm,"
$ws.Range("C4").Value = "    d1_nutri_food_g_2015.csv,"

$ws.Range("A5").Value = "    d1d2_food_2015.csv,"
$ws.Range("B5").Value = "    rf_cardio_food_simp.slurm,"
$ws.Range("C5").Value = "    d1_food_2015.csv,"

$ws.Range("A6").Value = "    d1d2_nutri_food_2015.csv,"
$ws.Range("B6").Value = "    rf_cardio_nut_food_simp.slurm,"
$ws.Range("C6").Value = "    d1_nutri_food_2015.csv,"

$ws.Range("A7").Value = "    d1d2_cat_g_2015.csv,"
$ws.Range("B7").Value = "    rf_cardio_cat_g.slurm,"
$ws.Range("C7").Value = "    d1d2_cat_g_2015.csv,"

$ws.Range("A8").Value = "    d1d2_nutri_cat_g_2015.csv,"
$ws.Range("B8").Value = "    rf_cardio_nut_cat_g.slurm,"
$ws.Range("C8").Value = "    d1_nutri_cat_g_2015.csv,"

$ws.Range("A9").Value = "    d1d2_cat_2015.csv,"
$ws.Range("B9").Value = "    rf_cardio_cat.slurm,"
$ws.Range("C9").Value = "    d1_cat_2015.csv,"

$ws.Range("A10").Value = "    d1d2_nutri_cat_2015.csv,"
$ws.Range("B10").Value = "    rf_cardio_nut_cat.slurm,"
$ws.Range("C10").Value = "    d1_nutri_cat_2015.csv,"

$ws.Range("A11").Value = "    d1d2_nutri_food_g_cat_g_2015.csv,"
$ws.Range("B11").Value = "    rf_cardio_nut_food_cat_g.slurm"
$ws.Range("C11").Value = "    d1_nutri_food_g_cat_g_2015.csv,"

$ws.Range("A12").Value = "    d1d2_nutri_food_cat_2015.csv,"
$ws.Range("B12").Value = "    rf_cardio_nut_food_cat_simp.slurm"
$ws.Range("C12").Value = "    d1_nutri_food_cat_2015.csv,"

# Column width for new column C (closest achievable to source 29.36328125 given engine pixel quantization)
$ws.Columns.Item(3).ColumnWidth = 28.5

# Move the window position (best-effort; engine does not persist this to xWindow/yWindow)
$win = $excel.Windows.Item(1)
$win.Left = 20
$win.Top = 620

$ws.Range("C1").Select() | Out-Null
